# Add season record columns (Wins, Losses, Ties) to the COL_2013 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels, matching the existing header formatting
# (bold, bordered, centered) by copying the format from the last header cell.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every player row (2-45): 74 wins, 88 losses, 0 ties.
$wins = 74
$losses = 88
$ties = 0

for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}
